# Arreglo de fallo encontrado en estimacion presupuesto y empiece de practica 8
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "I+d": update personnel-cost formulas (D5:D8) with the corrected
# multipliers discovered while fixing the budget estimation.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("I+d")
$ws1.Range("D5").Formula = "=(3*(2800*(1+B6)))*1"
$ws1.Range("D6").Formula = "=(5*(3050*(1+B6)))*5.4"
$ws1.Range("D7").Formula = "=(5*(2887*(1+B6)))*5.4"
$ws1.Range("D8").Formula = "=(3*(4200*(1+B6))*3)"

# ---------------------------------------------------------------------------
# Sheet "i+d+risks": insert a new row for executive-hiring costs above the
# existing "personal programador" row, then fix up the formulas for the
# rows that used to be D5:D8 (now D6:D9) with the same corrected
# multipliers, plus the new row D5.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("i+d+risks")
$ws2.Rows("5:5").Insert()

$ws2.Range("C5:D5").ClearFormats()
$ws2.Range("C5").Value = "Total gastos de contratación de ejecutivos"
$ws2.Range("D5").Formula = "=(4*(2800*(1+B7)))*5"

$ws2.Range("D6").Formula = "=4*(2800*(1+B7))"
$ws2.Range("D7").Formula = "=(7*(3050*(1+B7)))*5.4"
$ws2.Range("D8").Formula = "=(4*(2887*(1+B7)))*5.4"
$ws2.Range("D9").Formula = "=(6*(4200*(1+B7)))*3"
$ws2.Range("D4").Formula = "=SUM(D5:D9)"

# ---------------------------------------------------------------------------
# Sheet "Flujo de Cajas + VAN + TIR": revenue split now derives from the
# total expenses instead of fixed numbers, and the shared-formula anchor
# for row 6 moves from D6 to C6.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Flujo de Cajas + VAN + TIR")
$ws3.Range("F2").Formula = "=B18*0.4375"
$ws3.Range("J2").Formula = "=B18*0.5625"
$ws3.Range("D6:J6").Formula = "=`$C`$6"

$ws3.Range("B18").Value = $ws2.Range("D41").Value
$ws3.Range("B23").Value = $ws2.Range("D4").Value

$wb.RefreshAll()
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------------
# Restore the view state captured in the authored workbook.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 1
$ws1.Range("D8").Select()

$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollColumn = 6
$ws2.Range("L17").Select()

$ws3.Activate()
$ws3.Range("E19").Select()

$ws2.Activate()
